$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 3503
$ws.Range("I3").Value = 3643
$ws.Range("H4").Value = 1668
$ws.Range("I4").Value = 851
$ws.Range("I5").Value = 338
$ws.Range("I6").Value = 4079
$ws.Range("H7").Value = 25979
$ws.Range("I7").Value = 12414

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I3").Value = 31
$ws.Range("I6").Value = 52
$ws.Range("I7").Value = 140

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I2").Value = 38
$ws.Range("I7").Value = 141

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 127
$ws.Range("I6").Value = 110
$ws.Range("I7").Value = 393

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I3").Value = 83
$ws.Range("I4").Value = 16
$ws.Range("I7").Value = 230

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 170
$ws.Range("I5").Value = 11
$ws.Range("I6").Value = 160
$ws.Range("I7").Value = 480

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I6").Value = 36
$ws.Range("I7").Value = 116

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 90
$ws.Range("I7").Value = 274

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I6").Value = 86
$ws.Range("I7").Value = 404
$ws.Range("I8").Value = 759
$ws.Range("I9").Value = 57
$ws.Range("I11").Value = 198
$ws.Range("I13").Value = 20
$ws.Range("I19").Value = 329
$ws.Range("I20").Value = 306
$ws.Range("I24").Value = 33
$ws.Range("I29").Value = 807
$ws.Range("I31").Value = 116
$ws.Range("I33").Value = 560
$ws.Range("I34").Value = 57
$ws.Range("I36").Value = 171
$ws.Range("I37").Value = 393
$ws.Range("I42").Value = 428
$ws.Range("I47").Value = 84
$ws.Range("I48").Value = 162
$ws.Range("I49").Value = 101
$ws.Range("I52").Value = 271
$ws.Range("I54").Value = 277
$ws.Range("H63").Value = 208
$ws.Range("I63").Value = 48
$ws.Range("I65").Value = 274
$ws.Range("I67").Value = 480
$ws.Range("I71").Value = 36
$ws.Range("I76").Value = 189
$ws.Range("I77").Value = 68
$ws.Range("I78").Value = 176
$ws.Range("I79").Value = 322
$ws.Range("I81").Value = 12
$ws.Range("I83").Value = 248
$ws.Range("I85").Value = 574
$ws.Range("I88").Value = 113
$ws.Range("I89").Value = 140
$ws.Range("I90").Value = 155
$ws.Range("I91").Value = 151
$ws.Range("I93").Value = 67
$ws.Range("I94").Value = 111
$ws.Range("I96").Value = 141
$ws.Range("I98").Value = 79
$ws.Range("I99").Value = 230
$ws.Range("H101").Value = 25979
$ws.Range("I101").Value = 12414

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I3").Value = 94
$ws.Range("I6").Value = 46
$ws.Range("I7").Value = 248

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 132
$ws.Range("I3").Value = 201
$ws.Range("I7").Value = 560

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("I3").Value = 10
$ws.Range("I7").Value = 101

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I2").Value = 63
$ws.Range("I7").Value = 277

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 242
$ws.Range("I3").Value = 278
$ws.Range("I6").Value = 218
$ws.Range("I7").Value = 807

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 129
$ws.Range("I3").Value = 91
$ws.Range("I6").Value = 91
$ws.Range("I7").Value = 329

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I6").Value = 89
$ws.Range("I7").Value = 162

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I6").Value = 80
$ws.Range("I7").Value = 189

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I3").Value = 232
$ws.Range("I7").Value = 574

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I2").Value = 38
$ws.Range("I7").Value = 86

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I4").Value = 37
$ws.Range("I6").Value = 112
$ws.Range("I7").Value = 428

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range("I5").Value = 7
$ws.Range("I6").Value = 20

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I6").Value = 72
$ws.Range("I7").Value = 176

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I3").Value = 38
$ws.Range("I6").Value = 46

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("I2").Value = 10
$ws.Range("I7").Value = 33

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I3").Value = 52
$ws.Range("I6").Value = 48
$ws.Range("I7").Value = 151

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 93
$ws.Range("I7").Value = 322

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I6").Value = 99
$ws.Range("I7").Value = 306

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I3").Value = 52
$ws.Range("I4").Value = 7
$ws.Range("I7").Value = 171

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I2").Value = 19
$ws.Range("I7").Value = 67

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I6").Value = 64
$ws.Range("I7").Value = 271

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("I3").Value = 19
$ws.Range("I7").Value = 57

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I3").Value = 20
$ws.Range("I4").Value = 9
$ws.Range("I7").Value = 111

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I6").Value = 31
$ws.Range("I7").Value = 84

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I6").Value = 49
$ws.Range("I7").Value = 79

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I3").Value = 42
$ws.Range("I6").Value = 46
$ws.Range("I7").Value = 198

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("I3").Value = 21
$ws.Range("I7").Value = 57

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I3").Value = 42
$ws.Range("I7").Value = 113

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I3").Value = 211
$ws.Range("I7").Value = 759

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I4").Value = 14
$ws.Range("I7").Value = 155

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("I3").Value = 13
$ws.Range("I7").Value = 36

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("I2").Value = 17
$ws.Range("I7").Value = 68

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I4").Value = 21
$ws.Range("I7").Value = 404

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range("I3").Value = 3
$ws.Range("I6").Value = 12
